$d = $word.ActiveDocument

# Locate the paragraph that ends with "Worst game (most shots): 77"
# (the last stat line of the Expectimax section).
$found = $d.Content
$found.Find.Execute("Worst game (most shots): 77", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$targetEnd = $found.End

# Map that character position to its 1-based Paragraphs index so we can grab
# the Range of the paragraph that immediately follows it (the existing blank
# separator paragraph right before "Heatmap Strategy/Agent").
$paraIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.End -eq ($targetEnd + 1)) {
        $paraIndex = $i
        break
    }
}
if ($paraIndex -eq -1) {
    throw "Could not locate the 'Worst game (most shots): 77' paragraph"
}
$anchor = $d.Paragraphs.Item($paraIndex + 1).Range
$anchor.Collapse(1)

# Build the block of new paragraphs to splice in right before that existing
# blank separator: a new blank line, four new stat lines, and then a couple
# of trailing blank paragraphs so the existing blank separator before
# "Heatmap Strategy/Agent" is preserved after the splice.
$newParagraphsXml = ( `
    '<w:p/>' + `
    '<w:p><w:r><w:t>Average shots to win: 50.1</w:t></w:r></w:p>' + `
    '<w:p><w:r><w:t>Best game (fewest shots): 35</w:t></w:r></w:p>' + `
    '<w:p><w:r><w:t>Worst game (most shots): 68</w:t></w:r></w:p>' + `
    '<w:p><w:r><w:t>Last 5 games: [58, 44, 58, 53, 47]</w:t></w:r></w:p>' + `
    '<w:p/>' + `
    '<w:p/>' `
)

$packageXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
    '<pkg:xmlData>' + `
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
    '<w:body>' + $newParagraphsXml + '</w:body>' + `
    '</w:document>' + `
    '</pkg:xmlData>' + `
    '</pkg:part>' + `
    '</pkg:package>'

$anchor.InsertXML($packageXml)
